$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Paragraph 1 (Felhasználók): "(email, ..." -> "(,id,email, ..."
#   New text: ","(no underline) + "id"(underline) + ",email"(no underline)
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(1).Range
$r = $d.Range($p.Start, $p.End)
$r.Find.Execute("email", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
$insertPos = $r.Start
$r.InsertBefore(",id,email")
$idRange = $d.Range($insertPos + 1, $insertPos + 3)
$idRange.Underline = 1

# ---------------------------------------------------------------------------
# Paragraph 2 (Éttermek): "(email, ..." -> "(id,email, ..."
#   New text: "id"(underline) + ",email"(no underline)
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(2).Range
$r = $d.Range($p.Start, $p.End)
$r.Find.Execute("email", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
$insertPos = $r.Start
$r.InsertBefore("id,email")
$idRange = $d.Range($insertPos, $insertPos + 2)
$idRange.Underline = 1

# ---------------------------------------------------------------------------
# Paragraph 3 (Étlap): "éttermek_email," -> "éttermek_id," (stays underlined)
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(3).Range
$r = $d.Range($p.Start, $p.End)
$r.Find.Execute("éttermek_email,", $true, $false, $false, $false, $false, $true, 1, $false, "éttermek_id,", 2)

# ---------------------------------------------------------------------------
# Paragraph 4 (Nyitvatartás): "étterem_email" -> "étterem_id" (stays underlined)
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(4).Range
$r = $d.Range($p.Start, $p.End)
$r.Find.Execute("étterem_email", $true, $false, $false, $false, $false, $true, 1, $false, "étterem_id", 2)

# ---------------------------------------------------------------------------
# Paragraph 6 (Értékelés): "étterem_email, felh.email, " -> "étterem_id, felh.id, " (stays underlined)
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(6).Range
$r = $d.Range($p.Start, $p.End)
$r.Find.Execute("étterem_email, felh.email, ", $true, $false, $false, $false, $false, $true, 1, $false, "étterem_id, felh.id, ", 2)

# ---------------------------------------------------------------------------
# Paragraph 7 (Kedvenc): "étterem_email, felh.email" -> "étterem_id, felh.id" (stays underlined)
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(7).Range
$r = $d.Range($p.Start, $p.End)
$r.Find.Execute("étterem_email, felh.email", $true, $false, $false, $false, $false, $true, 1, $false, "étterem_id, felh.id", 2)

# ---------------------------------------------------------------------------
# Paragraph 9 (Helyfoglalás): "étterem_email, felh.email" -> "étterem_id, felh.id" (stays underlined)
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(9).Range
$r = $d.Range($p.Start, $p.End)
$r.Find.Execute("étterem_email, felh.email", $true, $false, $false, $false, $false, $true, 1, $false, "étterem_id, felh.id", 2)

# ---------------------------------------------------------------------------
# Paragraph 11 (Hibajelentés): "étterem_email, felh.email" -> "étterem_id, felh.id" (stays underlined)
# and the ",tipus,leírás" segment (plus the underlined "," before "tipus") moves
# from before the _GoBack bookmark to after it.
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(11).Range
$r = $d.Range($p.Start, $p.End)
$r.Find.Execute("étterem_email, felh.email", $true, $false, $false, $false, $false, $true, 1, $false, "étterem_id, felh.id", 2)

# Remove the "," (underlined) + "tipus" + ",leírás" run text that currently sits
# right before the bookmark, then re-insert the underlined "," after the bookmark.
$p = $d.Paragraphs.Item(11).Range
$r = $d.Range($p.Start, $p.End)
$r.Find.Execute(",tipus,leírás", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

$p = $d.Paragraphs.Item(11).Range
$r = $d.Range($p.Start, $p.End)
$r.Find.Execute(")", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
$insertPos = $r.Start
$r.InsertBefore(",tipus,leírás)")
$commaRange = $d.Range($insertPos, $insertPos + 1)
$commaRange.Underline = 1
